# Updated cryptos list on Thu Mar  7 06:16:23 UTC 2024 with GitHub Actions
# Refresh Price (col D) and Volume(1h) (col E) figures for each coin row,
# and for the few rows whose rank order changed (36/37, 40/41, 44/45)
# swap the Coin name + Link (cols B/C) to match the new ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading "'" on numeric-looking Price strings forces Excel to keep them
# as text (preserving trailing zeros, e.g. "1.00", "0.140") instead of
# silently coercing to a Double.

$ws.Range("D2").Value = '65.855.05'
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("D3").Value = '3.751.55'
$ws.Range("E3").Value = '  -0.87%  '
$ws.Range("D4").Value = '''0.998'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''424.66'
$ws.Range("E5").Value = '  +4.63%  '
$ws.Range("D6").Value = '''136.77'
$ws.Range("E6").Value = '  +3.40%  '
$ws.Range("D7").Value = '''0.619'
$ws.Range("E7").Value = '  +2.04%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -0.25%  '
$ws.Range("D10").Value = '''0.149'
$ws.Range("E10").Value = '  -11.36%  '
$ws.Range("D11").Value = '''0.0000298'
$ws.Range("E11").Value = '  -17.47%  '
$ws.Range("D12").Value = '''42.01'
$ws.Range("E12").Value = '  +3.93%  '
$ws.Range("D13").Value = '''10.29'
$ws.Range("E13").Value = '  +6.07%  '
$ws.Range("D14").Value = '4.355.58'
$ws.Range("E14").Value = '  +0.01%  '
$ws.Range("D15").Value = '''14.94'
$ws.Range("E15").Value = '  +3.89%  '
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("D17").Value = '3.743.10'
$ws.Range("E17").Value = '  -0.67%  '
$ws.Range("D18").Value = '''19.68'
$ws.Range("E18").Value = '  +1.16%  '
$ws.Range("D19").Value = '''1.11'
$ws.Range("E19").Value = '  +4.44%  '
$ws.Range("D20").Value = '65.804.01'
$ws.Range("E20").Value = '  -1.18%  '
$ws.Range("D21").Value = '''400.09'
$ws.Range("E21").Value = '  -3.66%  '
$ws.Range("D22").Value = '''14.80'
$ws.Range("E22").Value = '  +2.80%  '
$ws.Range("D23").Value = '''3.17'
$ws.Range("E23").Value = '  +4.96%  '
$ws.Range("D24").Value = '''83.82'
$ws.Range("E24").Value = '  -1.50%  '
$ws.Range("D25").Value = '''36.17'
$ws.Range("E25").Value = '  -0.57%  '
$ws.Range("D26").Value = '''9.66'
$ws.Range("E26").Value = '  +34.06%  '
$ws.Range("E27").Value = '  +4.46%  '
$ws.Range("D28").Value = '''9.74'
$ws.Range("E28").Value = '  +4.83%  '
$ws.Range("D29").Value = '''5.42'
$ws.Range("E29").Value = '  -4.22%  '
$ws.Range("D30").Value = '''13.59'
$ws.Range("E30").Value = '  +9.98%  '
$ws.Range("D31").Value = '''696.86'
$ws.Range("E31").Value = '  +0.60%  '
$ws.Range("D32").Value = '''0.129'
$ws.Range("E32").Value = '  +10.36%  '
$ws.Range("D33").Value = '''2.76'
$ws.Range("E33").Value = '  +1.28%  '
$ws.Range("D34").Value = '''40.04'
$ws.Range("E34").Value = '  +2.78%  '
$ws.Range("D35").Value = '''0.998'
$ws.Range("E35").Value = '  -0.19%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '''0.147'
$ws.Range("E36").Value = '  -4.86%  '
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").Value = '''5.55'
$ws.Range("E37").Value = '  +31.19%  '
$ws.Range("D38").Value = '''56.06'
$ws.Range("E38").Value = '  +1.66%  '
$ws.Range("D39").Value = '''0.0465'
$ws.Range("E39").Value = '  +2.05%  '
$ws.Range("B40").Value = 'ThetaToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D40").Value = '''2.93'
$ws.Range("E40").Value = '  +2.53%  '
$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").Value = '''2.63'
$ws.Range("E41").Value = '  +36.67%  '
$ws.Range("E42").Value = '  +0.69%  '
$ws.Range("D43").Value = '''0.140'
$ws.Range("E43").Value = '  +3.88%  '
$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").Value = '''3.22'
$ws.Range("E44").Value = '  +2.25%  '
$ws.Range("B45").Value = 'PEPE'
$ws.Range("C45").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D45").Value = '0.0₃0653'
$ws.Range("E45").Value = '  -12.59%  '
$ws.Range("D46").Value = '''3.32'
$ws.Range("E46").Value = '  +2.73%  '
$ws.Range("D47").Value = '''0.318'
$ws.Range("E47").Value = '  +9.49%  '
$ws.Range("D48").Value = '''2.66'
$ws.Range("E48").Value = '  +4.10%  '
$ws.Range("D49").Value = '''2.03'
$ws.Range("E49").Value = '  -0.56%  '
$ws.Range("D50").Value = '''138.92'
$ws.Range("E50").Value = '  -3.76%  '
$ws.Range("D51").Value = '''2.75'
$ws.Range("E51").Value = '  -2.47%  '
